$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 height changed (90 -> 45)
$ws.Rows.Item(3).RowHeight = 45

# New row 10 content
$ws.Range("D10").Value = "Danh sách tài khoản"
$ws.Range("E10").Value = "`n"
$ws.Range("F10").Value = "SELECT"
$ws.Range("G10").Value = "Nhanvien,`nTaikhoan`nTaikhoan_Nhanvien`nNhomquyen"
$ws.Range("J10").Value = "NhanvienID,Tennhanvien,tentaikhoan,Tenquyen`n (multi record)"

$ws.Rows.Item(10).RowHeight = 60

# Borders for D10 and F10 (left+right thin only)
$ws.Range("D10").Borders.Item(7).LineStyle = 1
$ws.Range("D10").Borders.Item(10).LineStyle = 1
$ws.Range("F10").Borders.Item(7).LineStyle = 1
$ws.Range("F10").Borders.Item(10).LineStyle = 1

# Border for J10 (left+right thin) + wrap text
$ws.Range("J10").Borders.Item(7).LineStyle = 1
$ws.Range("J10").Borders.Item(10).LineStyle = 1
$ws.Range("J10").WrapText = $true

# Vertical alignment top for D10/F10/J10
$ws.Range("D10").VerticalAlignment = -4160
$ws.Range("F10").VerticalAlignment = -4160
$ws.Range("J10").VerticalAlignment = -4160

# Wrap text for E10/G10 (no border, no vertical alignment)
$ws.Range("E10").WrapText = $true
$ws.Range("G10").WrapText = $true

# Sheet view changes
$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("K10").Select()
